$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Thursday")

$rows = @(
    @("Mealer", "Brock", "rnc.mealer.txt"),
    @("Blackburn", "Marsha", "rnc.blackburn.txt"),
    @("Fallin", "Mary", "rnc.fallin.txt"),
    @("Shin", "Lisa", "rnc.shin.txt"),
    @("Falwell", "Jerry", "rnc.falwell.txt"),
    @("Thiel", "Peter", "rnc.thiel.txt"),
    @("Burns", "Mark", "rnc.burns.txt"),
    @("Barrack", "Tom", "rnc.barrack.txt"),
    @("Tarkenton", "Fran", "rnc.tarkenton.txt"),
    @("Perkins", "Tony", "rnc.perkins.txt")
)

$ws.Range("F11").ClearFormats()

$r = 6
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = "Thursday"
    $ws.Cells.Item($r, 5).Value = "speech"
    $ws.Cells.Item($r, 6).Value = "Cision"
    $r = $r + 1
}

$ws.Range("F16").Select()
$ws.Activate()
